# Auto-generated edit script: update Sheets via scheduled runner
# Applies updated market-board price data to leve-profit rows across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()

$ws.Range("H80").Value = 288572.34
$ws.Range("I80").Value = 591.2
$ws.Range("J80").Value = 648548.75
$ws.Range("K80").Value = 1773.6
$ws.Range("L80").Value = 1945646.25
$ws.Range("M80").Value = -775.6000000000001
$ws.Range("N80").Value = -1947642.25

$ws.Range("H83").Value = 288572.34
$ws.Range("I83").Value = 591.2
$ws.Range("J83").Value = 648548.75
$ws.Range("K83").Value = 5320.8
$ws.Range("L83").Value = 5836938.75
$ws.Range("M83").Value = -328.8000000000002
$ws.Range("N83").Value = -5846922.75

$ws.Range("H98").Value = 717
$ws.Range("I98").Value = 694.5
$ws.Range("J98").Value = 867
$ws.Range("K98").Value = 694.5
$ws.Range("L98").Value = 867
$ws.Range("M98").Value = 803.5
$ws.Range("N98").Value = -3863

$ws.Range("H113").Value = 2747.8948
$ws.Range("I113").Value = 2639.2307
$ws.Range("J113").Value = 2983.3333
$ws.Range("K113").Value = 2639.2307
$ws.Range("L113").Value = 2983.3333
$ws.Range("M113").Value = 614.7692999999999
$ws.Range("N113").Value = -9491.3333

$ws.Range("H121").Value = 1058.4286
$ws.Range("J121").Value = 1161.8
$ws.Range("L121").Value = 3485.4
$ws.Range("N121").Value = -6979.4

$ws.Range("H122").Value = 717
$ws.Range("I122").Value = 694.5
$ws.Range("J122").Value = 867
$ws.Range("K122").Value = 2083.5
$ws.Range("L122").Value = 2601
$ws.Range("M122").Value = 366.5
$ws.Range("N122").Value = -7501

$ws.Range("H141").Value = 1533.7142
$ws.Range("I141").Value = 1533.7142
$ws.Range("K141").Value = 4601.142599999999
$ws.Range("M141").Value = 578.8574000000008

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8034.3936
$ws.Range("I32").Value = 4895.4644
$ws.Range("J32").Value = 34401.4
$ws.Range("K32").Value = 4895.4644
$ws.Range("L32").Value = 34401.4
$ws.Range("M32").Value = -4608.4644
$ws.Range("N32").Value = -34975.4

$ws.Range("H61").Value = 1546.2222
$ws.Range("I61").Value = 1489.5
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1489.5
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1277.5
$ws.Range("N61").Value = -2424

$ws.Range("H102").Value = 2373.8096
$ws.Range("I102").Value = 2396.7896
$ws.Range("J102").Value = 2155.5
$ws.Range("K102").Value = 2396.7896
$ws.Range("L102").Value = 2155.5
$ws.Range("M102").Value = -774.7896000000001
$ws.Range("N102").Value = -5399.5

$ws.Range("H136").Value = 1546.2222
$ws.Range("I136").Value = 1489.5
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 4468.5
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -1918.5
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 23257812
$ws.Range("I86").Value = 31251702
$ws.Range("J86").Value = 2859.5454
$ws.Range("K86").Value = 31251702
$ws.Range("L86").Value = 2859.5454
$ws.Range("M86").Value = -31250579
$ws.Range("N86").Value = -5105.5454

$ws.Range("H89").Value = 23257812
$ws.Range("I89").Value = 31251702
$ws.Range("J89").Value = 2859.5454
$ws.Range("K89").Value = 156258510
$ws.Range("L89").Value = 14297.727
$ws.Range("M89").Value = -156252894
$ws.Range("N89").Value = -25529.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 5048.3335
$ws.Range("I70").Value = 4672.5
$ws.Range("J70").Value = 5800
$ws.Range("K70").Value = 14017.5
$ws.Range("L70").Value = 17400
$ws.Range("M70").Value = -13702.5
$ws.Range("N70").Value = -18030

$ws.Range("H73").Value = 5048.3335
$ws.Range("I73").Value = 4672.5
$ws.Range("J73").Value = 5800
$ws.Range("K73").Value = 14017.5
$ws.Range("L73").Value = 17400
$ws.Range("M73").Value = -12925.5
$ws.Range("N73").Value = -19584

$ws.Range("H113").Value = 563.3333
$ws.Range("I113").Value = 1222
$ws.Range("J113").Value = 520.8387
$ws.Range("K113").Value = 3666
$ws.Range("L113").Value = 1562.5161
$ws.Range("M113").Value = -1496
$ws.Range("N113").Value = -5902.5161

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4124.885
$ws.Range("I80").Value = 4115.8823
$ws.Range("J80").Value = 4141.8887
$ws.Range("K80").Value = 4115.8823
$ws.Range("L80").Value = 4141.8887
$ws.Range("M80").Value = -3117.8823
$ws.Range("N80").Value = -6137.8887

$ws.Range("H83").Value = 4124.885
$ws.Range("I83").Value = 4115.8823
$ws.Range("J83").Value = 4141.8887
$ws.Range("K83").Value = 20579.4115
$ws.Range("L83").Value = 20709.4435
$ws.Range("M83").Value = -15587.4115
$ws.Range("N83").Value = -30693.4435

$ws.Range("H126").Value = 3940.2083
$ws.Range("I126").Value = 3458.0588
$ws.Range("J126").Value = 5111.143
$ws.Range("K126").Value = 10374.1764
$ws.Range("L126").Value = 15333.429
$ws.Range("M126").Value = -7904.1764
$ws.Range("N126").Value = -20273.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3480.4614
$ws.Range("I68").Value = 3382.5557
$ws.Range("K68").Value = 3382.5557
$ws.Range("M68").Value = -2633.5557

$ws.Range("H71").Value = 3480.4614
$ws.Range("I71").Value = 3382.5557
$ws.Range("K71").Value = 16912.7785
$ws.Range("M71").Value = -13168.7785

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 13725.667
$ws.Range("J41").Value = 13725.667
$ws.Range("L41").Value = 13725.667
$ws.Range("N41").Value = -14505.667

$ws.Range("H136").Value = 340.47223
$ws.Range("I136").Value = 342.89655
$ws.Range("J136").Value = 330.42856
$ws.Range("K136").Value = 1028.68965
$ws.Range("L136").Value = 991.28568
$ws.Range("M136").Value = 1521.31035
$ws.Range("N136").Value = -6091.28568
